# Create a new "Table" table style (as introduced by the commit
# "Create Table style for tables"):
#
#   <w:style w:type="table" w:customStyle="1" w:styleId="Table">
#     <w:name w:val="Table"/>
#     <w:basedOn w:val="TableNormal"/>
#     <w:uiPriority w:val="99"/>
#     <w:rsid w:val="00296345"/>
#     <w:pPr>
#       <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
#     </w:pPr>
#     <w:tblPr>
#       <w:tblBorders>
#         <w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/>
#         <w:left w:val="single" w:sz="4" w:space="0" w:color="auto"/>
#         <w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/>
#         <w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/>
#         <w:insideH w:val="single" w:sz="4" w:space="0" w:color="auto"/>
#         <w:insideV w:val="single" w:sz="4" w:space="0" w:color="auto"/>
#       </w:tblBorders>
#     </w:tblPr>
#   </w:style>

$d = $word.ActiveDocument

# wdStyleTypeTable = 3
$tableStyle = $d.Styles.Add("Table", 3)

# <w:basedOn w:val="TableNormal"/>
$tableStyle.BaseStyle = "TableNormal"

# <w:uiPriority w:val="99"/>
$tableStyle.Priority = 99

# <w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>
$pf = $tableStyle.ParagraphFormat
$pf.SpaceAfter = 0
# wdLineSpaceSingle = 0 -> w:line="240" w:lineRule="auto"
$pf.LineSpacingRule = 0
